# Reference Genome Assembly Test suite
#
# The "startup" sheet lists the tabs that the Commons/Neo4j vs Web
# comparison runs against. The first row used to document the "Cases"
# tab; this suite instead documents the "Participants" tab, so rename
# the label and refresh the sheet's view/layout metrics that Excel
# recalculates as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# 1. Rename the tab label in A2: CasesTab -> ParticipantsTab
$ws.Range("A2").Value = "ParticipantsTab"

# 2. Column A now holds a longer label ("ParticipantsTab"); widen it and
#    drop the stale "best fit" sizing in favor of an explicit width.
$ws.Columns.Item(1).ColumnWidth = 15.7

# 3. Refresh the (wrapped-text) row heights for rows 2-4 to their
#    recalculated values.
$ws.Rows.Item(2).RowHeight = 165
$ws.Rows.Item(3).RowHeight = 180
$ws.Rows.Item(4).RowHeight = 210

# 4. Reset the view: scroll back to the top-left corner (clearing the
#    old topLeftCell="A2") and select the entire sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
[void]$ws.Cells.Select()
